$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New GitHub Admin Log row recorded for the delete-team action.
$ws.Range("A3").Value = "2025-07-23 12:38:40"
$ws.Range("B3").Value = "delete-team"
$ws.Range("C3").Value = "new-organization97"
$ws.Range("D3").Value = "firstteam"
$ws.Range("E3").Value = "task-repo"

# Column I holds "True"/"False" as literal text (see I2), so force text
# entry with a leading apostrophe rather than letting Excel's type
# inference turn it into a Boolean. Reset the cell style back to Normal
# afterwards so the quote-prefix marker doesn't linger as cell formatting.
$ws.Range("I3").Value = "'False"
$ws.Range("I3").Style = "Normal"
